# Applies the "initial version of stimuli order" edit: re-shuffles the
# stimulus rows (columns B-E) of the 18_cues sequence sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0, 55, "face/face085.png", "scheitern", "face"),
    @(1, 124, "face/face121.png", "rasen", "face"),
    @(2, 100, "face/face099.png", "wenden", "face"),
    @(3, 73, "car/car069.png", "segeln", "car"),
    @(4, 3, "face/face088.png", "fesseln", "face"),
    @(5, 79, "face/face101.png", "hauen", "face"),
    @(6, 81, "face/face097.png", "liefern", "face"),
    @(7, 7, "car/car068.png", "kehren", "car"),
    @(8, 35, "car/car072.png", "laufen", "car"),
    @(9, 42, "car/car077.png", "biegen", "car"),
    @(10, 14, "car/car098.png", "starten", "car"),
    @(11, 60, "car/car096.png", "klappen", "car"),
    @(12, 120, "face/face109.png", "drehen", "face"),
    @(13, 61, "car/car070.png", "schenken", "car"),
    @(14, 126, "face/face108.png", "langen", "face"),
    @(15, 97, "car/car116.png", "opfern", "car"),
    @(16, 123, "car/car073.png", "jubeln", "car"),
    @(17, 43, "face/face091.png", "haken", "face"),
    @(18, 31, "face/face069.png", "schmecken", "face"),
    @(19, 66, "car/car081.png", "pflegen", "car"),
    @(20, 40, "face/face104.png", "kaufen", "face"),
    @(21, 98, "face/face093.png", "mieten", "face"),
    @(22, 99, "car/car074.png", "tagen", "car"),
    @(23, 46, "car/car087.png", "bitten", "car"),
    @(24, 85, "face/face064.png", "drohen", "face"),
    @(25, 24, "face/face070.png", "sondern", "face"),
    @(26, 5, "car/car091.png", "hupen", "car"),
    @(27, 11, "car/car088.png", "husten", "car"),
    @(28, 64, "face/face090.png", "saufen", "face"),
    @(29, 125, "car/car113.png", "schicken", "car"),
    @(30, 105, "face/face086.png", "hoffen", "face"),
    @(31, 94, "car/car095.png", "ehren", "car")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $excelRow = $i + 2
    $ws.Cells.Item($excelRow, 1).Value = $row[0]
    $ws.Cells.Item($excelRow, 2).Value = $row[1]
    $ws.Cells.Item($excelRow, 3).Value = $row[2]
    $ws.Cells.Item($excelRow, 4).Value = $row[3]
    $ws.Cells.Item($excelRow, 5).Value = $row[4]
}
